$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 334
$ws.Range("B334").Value = 6867195
$ws.Range("F334").Value = "SzegedCsanad Grosics"
$ws.Range("G334").Value = "Pecsi MFC"
$ws.Range("K334").Value = 2.1
$ws.Range("L334").Value = 3.4
$ws.Range("M334").Value = 2.9
$ws.Range("N334").Value = 1.95
$ws.Range("O334").Value = 3.4
$ws.Range("P334").Value = 4
$ws.Range("Q334").Value = -0.5
$ws.Range("R334").Value = 1.975
$ws.Range("S334").Value = 1.825
$ws.Range("T334").Value = 1.75
$ws.Range("U334").Value = 1.775
$ws.Range("V334").Value = 2.025
$ws.Range("W334").Value = 0.95
$ws.Range("Z334").Value = 0.9750000000000001
$ws.Range("AA334").Value = -1
$ws.Range("AC334").Value = 1.025

# Row 335
$ws.Range("B335").Value = 6867248
$ws.Range("F335").Value = "Gyori ETO"
$ws.Range("G335").Value = "FC Ajka"
$ws.Range("I335").Value = 0
$ws.Range("J335").Value = "H"
$ws.Range("K335").Value = 1.5
$ws.Range("L335").Value = 4
$ws.Range("M335").Value = 5
$ws.Range("N335").Value = 1.533
$ws.Range("O335").Value = 4.2
$ws.Range("P335").Value = 6
$ws.Range("Q335").Value = -1
$ws.Range("R335").Value = 1.9
$ws.Range("S335").Value = 1.9
$ws.Range("T335").Value = 2.5
$ws.Range("U335").Value = 1.925
$ws.Range("V335").Value = 1.875
$ws.Range("W335").Value = 0.5329999999999999
$ws.Range("X335").Value = -1
$ws.Range("Z335").Value = 0
$ws.Range("AA335").Value = -0
$ws.Range("AB335").Value = -1
$ws.Range("AC335").Value = 0.875

# Row 336
$ws.Range("B336").Value = 6867976
$ws.Range("F336").Value = "Tiszakecske FC"
$ws.Range("G336").Value = "Csakvari Tk"
$ws.Range("K336").Value = 2.5
$ws.Range("L336").Value = 3.2
$ws.Range("M336").Value = 2.5
$ws.Range("N336").Value = 2.375
$ws.Range("O336").Value = 3.2
$ws.Range("P336").Value = 3.2
$ws.Range("Q336").Value = -0.25
$ws.Range("R336").Value = 2
$ws.Range("S336").Value = 1.8
$ws.Range("T336").Value = 2
$ws.Range("U336").Value = 1.825
$ws.Range("V336").Value = 1.975
$ws.Range("X336").Value = 2.2
$ws.Range("Z336").Value = -0.5
$ws.Range("AA336").Value = 0.4
$ws.Range("AB336").Value = 0
$ws.Range("AC336").Value = -0

# Row 337
$ws.Range("B337").Value = 6871116
$ws.Range("F337").Value = "Budapest Honved"
$ws.Range("G337").Value = "BVSC Zuglo"
$ws.Range("I337").Value = 1
$ws.Range("J337").Value = "D"
$ws.Range("K337").Value = 1.85
$ws.Range("L337").Value = 3.5
$ws.Range("M337").Value = 3.5
$ws.Range("N337").Value = 1.75
$ws.Range("O337").Value = 3.6
$ws.Range("P337").Value = 4.75
$ws.Range("R337").Value = 1.725
$ws.Range("S337").Value = 2.075
$ws.Range("T337").Value = 2.25
$ws.Range("U337").Value = 1.925
$ws.Range("V337").Value = 1.875
$ws.Range("W337").Value = -1
$ws.Range("X337").Value = 2.6
$ws.Range("Z337").Value = -1
$ws.Range("AA337").Value = 1.075
$ws.Range("AB337").Value = -0.5
$ws.Range("AC337").Value = 0.4375

# Row 357
$ws.Range("B357").Value = 6867602
$ws.Range("F357").Value = "MTE 1904"
$ws.Range("G357").Value = "Vasas SC"
$ws.Range("H357").Value = 1
$ws.Range("I357").Value = 3
$ws.Range("J357").Value = "A"
$ws.Range("K357").Value = 5
$ws.Range("L357").Value = 3.75
$ws.Range("M357").Value = 1.55
$ws.Range("N357").Value = 6.5
$ws.Range("O357").Value = 4
$ws.Range("P357").Value = 1.4
$ws.Range("Q357").Value = 1.25
$ws.Range("R357").Value = 1.85
$ws.Range("S357").Value = 1.95
$ws.Range("T357").Value = 2.5
$ws.Range("U357").Value = 1.825
$ws.Range("V357").Value = 1.975
$ws.Range("W357").Value = -1
$ws.Range("Y357").Value = 0.3999999999999999
$ws.Range("Z357").Value = -1
$ws.Range("AA357").Value = 0.95
$ws.Range("AB357").Value = 0.825

# Row 358
$ws.Range("B358").Value = 6867554
$ws.Range("F358").Value = "Soroksar"
$ws.Range("G358").Value = "BVSC Zuglo"
$ws.Range("H358").Value = 2
$ws.Range("I358").Value = 1
$ws.Range("J358").Value = "H"
$ws.Range("K358").Value = 1.95
$ws.Range("L358").Value = 3.4
$ws.Range("M358").Value = 3.25
$ws.Range("N358").Value = 2.1
$ws.Range("O358").Value = 3.2
$ws.Range("P358").Value = 3.1
$ws.Range("Q358").Value = -0.25
$ws.Range("R358").Value = 1.9
$ws.Range("S358").Value = 1.9
$ws.Range("T358").Value = 2
$ws.Range("U358").Value = 1.9
$ws.Range("V358").Value = 1.9
$ws.Range("W358").Value = 1.1
$ws.Range("Y358").Value = -1
$ws.Range("Z358").Value = 0.8999999999999999
$ws.Range("AA358").Value = -1
$ws.Range("AB358").Value = 0.8999999999999999
